$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.507621333333333
$ws.Range("H2").Value = 7.522864
$ws.Range("I2").Value = 0.9939780200440224
$ws.Range("J2").Value = 0.9939780200440224
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.660964666666667
$ws.Range("N2").Value = 4.982894
$ws.Range("O2").Value = 0.09230299804773986
$ws.Range("P2").Value = 0.09230299804773986
$ws.Range("Q2").Value = 4.165070432046222
$ws.Range("R2").Value = 37.485633888416
$ws.Range("S2").Value = 0.09174715124361973
$ws.Range("T2").Value = 0.09174715124361973

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.507621333333333
$ws.Range("H3").Value = 7.522864
$ws.Range("I3").Value = 0.9939780200440224
$ws.Range("J3").Value = 0.9939780200440224
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.723353333333334
$ws.Range("N3").Value = 26.17006
$ws.Range("O3").Value = 0.4847735065384163
$ws.Range("P3").Value = 0.4847735065384162
$ws.Range("Q3").Value = 21.87486691687111
$ws.Range("R3").Value = 196.87380225184
$ws.Range("S3").Value = 0.481854210198853
$ws.Range("T3").Value = 0.4818542101988529

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.507621333333333
$ws.Range("H4").Value = 7.522864
$ws.Range("I4").Value = 0.9939780200440224
$ws.Range("J4").Value = 0.9939780200440224
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.610381
$ws.Range("N4").Value = 22.831143
$ws.Range("O4").Value = 0.4229234954138438
$ws.Range("P4").Value = 0.4229234954138438
$ws.Range("Q4").Value = 19.08395375039467
$ws.Range("R4").Value = 171.755583753552
$ws.Range("S4").Value = 0.4203766586015497
$ws.Range("T4").Value = 0.4203766586015496

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01519233333333333
$ws.Range("H5").Value = 0.045577
$ws.Range("I5").Value = 0.006021979955977724
$ws.Range("J5").Value = 0.006021979955977723
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.660964666666667
$ws.Range("N5").Value = 4.982894
$ws.Range("O5").Value = 0.09230299804773986
$ws.Range("P5").Value = 0.09230299804773986
$ws.Range("Q5").Value = 0.02523392887088889
$ws.Range("R5").Value = 0.227105359838
$ws.Range("S5").Value = 0.0005558468041201405
$ws.Range("T5").Value = 0.0005558468041201404

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Lgr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01519233333333333
$ws.Range("H6").Value = 0.045577
$ws.Range("I6").Value = 0.006021979955977724
$ws.Range("J6").Value = 0.006021979955977723
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.723353333333334
$ws.Range("N6").Value = 26.17006
$ws.Range("O6").Value = 0.4847735065384163
$ws.Range("P6").Value = 0.4847735065384162
$ws.Range("Q6").Value = 0.1325280916244445
$ws.Range("R6").Value = 1.19275282462
$ws.Range("S6").Value = 0.00291929633956338
$ws.Range("T6").Value = 0.002919296339563378

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Lgr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01519233333333333
$ws.Range("H7").Value = 0.045577
$ws.Range("I7").Value = 0.006021979955977724
$ws.Range("J7").Value = 0.006021979955977723
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.610381
$ws.Range("N7").Value = 22.831143
$ws.Range("O7").Value = 0.4229234954138438
$ws.Range("P7").Value = 0.4229234954138438
$ws.Range("Q7").Value = 0.1156194449456667
$ws.Range("R7").Value = 1.040575004511
$ws.Range("S7").Value = 0.002546836812294204
$ws.Range("T7").Value = 0.002546836812294204
